$wb = $excel.ActiveWorkbook

# --- Description sheet: insert a new row 5 ("Desc Body" / NPA description text) ---
# Copy row 1's formatting (Benchmark / long-title row) down onto the freshly
# inserted row 5 so the new cells inherit the sheet's existing label/value
# styling, then overwrite with the real content.
$wsDesc = $wb.Worksheets.Item("Description")
$wsDesc.Rows.Item(1).Copy()
$wsDesc.Rows.Item(5).Insert()

$wsDesc.Range("A5").Value = "Desc Body"
$wsDesc.Range("B5").Value = "The Commonwealth committed `$1.75 billion over five years for a National Partnership Agreement on Skills Reform (NPASR) to reform the vocational education and training system. The NPASR commenced in 2012 and expires in June 2017."
$wsDesc.Range("B5").Font.Color = 0
$wsDesc.Range("B5").WrapText = $true
$wsDesc.Rows.Item(5).RowHeight = 23.95

# --- Selection / active-sheet bookkeeping to match the authored workbook state ---
$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("B15").Select() | Out-Null

$wsDesc.Activate() | Out-Null
$wsDesc.Range("B9").Select() | Out-Null
